$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two recomputed error values (D22, C23) per the new sampling run
$ws.Range("D22").Value = 0.7115302104241067
$ws.Range("C23").Value = 0.3054124294241067

# Append the new ifoCAST-sampling diff row (row 24): copy the date-label
# formatting from the cell above (A23, style s="1") so the new label keeps
# the same bold/border/centered look as the rest of column A, then fill in
# the new values.
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A24").Value = "2025-09-04_diff"
$ws.Range("B24").Value = 0.05603945542410671
